# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the Malboro_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the authoritative XML diff (chore: update Sheets via scheduled runner).

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 458.29413
$ws.Range("I5").Value = 181.16667
$ws.Range("J5").Value = 609.4545000000001
$ws.Range("K5").Value = 181.16667
$ws.Range("L5").Value = 609.4545000000001
$ws.Range("M5").Value = -66.16667000000001
$ws.Range("N5").Value = -839.4545000000001
$ws.Range("H9").Value = 499.75
$ws.Range("I9").Value = 778.2857
$ws.Range("J9").Value = 109.8
$ws.Range("K9").Value = 778.2857
$ws.Range("L9").Value = 109.8
$ws.Range("M9").Value = -609.2857
$ws.Range("N9").Value = -447.8
$ws.Range("H51").Value = 6000
$ws.Range("J51").Value = 7000
$ws.Range("L51").Value = 7000
$ws.Range("N51").Value = -7968
$ws.Range("H53").Value = 5665.1816
$ws.Range("I53").Value = 859.9
$ws.Range("J53").Value = 9669.583000000001
$ws.Range("K53").Value = 859.9
$ws.Range("L53").Value = 9669.583000000001
$ws.Range("M53").Value = -222.9
$ws.Range("N53").Value = -10943.583
$ws.Range("H112").Value = 4420.794
$ws.Range("J112").Value = 2528.2334
$ws.Range("L112").Value = 7584.7002
$ws.Range("N112").Value = -9800.700199999999
$ws.Range("H116").Value = 4289.8657
$ws.Range("I116").Value = 3997.9365
$ws.Range("K116").Value = 3997.9365
$ws.Range("M116").Value = -555.9364999999998
$ws.Range("H125").Value = 4218.524
$ws.Range("I125").Value = 3730.6667
$ws.Range("J125").Value = 4869
$ws.Range("K125").Value = 33576.0003
$ws.Range("L125").Value = 43821
$ws.Range("M125").Value = -31116.0003
$ws.Range("N125").Value = -48741
$ws.Range("H132").Value = 4948.975
$ws.Range("I132").Value = 3625.5066
$ws.Range("K132").Value = 10876.5198
$ws.Range("M132").Value = -8346.5198
$ws.Range("H138").Value = 4749.74
$ws.Range("I138").Value = 2931.2903
$ws.Range("J138").Value = 5566.7246
$ws.Range("K138").Value = 8793.8709
$ws.Range("L138").Value = 16700.1738
$ws.Range("M138").Value = -3653.8709
$ws.Range("N138").Value = -26980.1738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1221578.5
$ws.Range("I32").Value = 1929725
$ws.Range("K32").Value = 1929725
$ws.Range("M32").Value = -1929438
$ws.Range("H45").Value = 4087.4
$ws.Range("I45").Value = 4570.467
$ws.Range("K45").Value = 4570.467
$ws.Range("M45").Value = -4193.467
$ws.Range("H61").Value = 15216.869
$ws.Range("I61").Value = 6503.8213
$ws.Range("K61").Value = 6503.8213
$ws.Range("M61").Value = -6291.8213
$ws.Range("H97").Value = 1433.3846
$ws.Range("I97").Value = 1292
$ws.Range("K97").Value = 1292
$ws.Range("M97").Value = -796
$ws.Range("H132").Value = 6730.1577
$ws.Range("J132").Value = 18436.334
$ws.Range("L132").Value = 55309.00199999999
$ws.Range("N132").Value = -60369.00199999999
$ws.Range("H133").Value = 97199.336
$ws.Range("J133").Value = 97199.336
$ws.Range("L133").Value = 97199.336
$ws.Range("N133").Value = -102259.336
$ws.Range("H135").Value = 265213
$ws.Range("J135").Value = 265213
$ws.Range("L135").Value = 265213
$ws.Range("N135").Value = -275353
$ws.Range("H136").Value = 15216.869
$ws.Range("I136").Value = 6503.8213
$ws.Range("K136").Value = 19511.4639
$ws.Range("M136").Value = -16961.4639

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17802.783
$ws.Range("I20").Value = 4310.269
$ws.Range("J20").Value = 35343.05
$ws.Range("K20").Value = 4310.269
$ws.Range("L20").Value = 35343.05
$ws.Range("M20").Value = -4063.269
$ws.Range("N20").Value = -35837.05
$ws.Range("H22").Value = 1038.4445
$ws.Range("I22").Value = 1038.4445
$ws.Range("K22").Value = 1038.4445
$ws.Range("M22").Value = -865.4445000000001
$ws.Range("H105").Value = 2167.842
$ws.Range("I105").Value = 2246.4119
$ws.Range("K105").Value = 2246.4119
$ws.Range("M105").Value = -499.4119000000001
$ws.Range("H134").Value = 8332.349
$ws.Range("I134").Value = 2702.1614
$ws.Range("K134").Value = 8106.4842
$ws.Range("M134").Value = -5571.4842
$ws.Range("H135").Value = 92423.336
$ws.Range("J135").Value = 92423.336
$ws.Range("L135").Value = 92423.336
$ws.Range("N135").Value = -102563.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10470
$ws.Range("H58").Value = 9516.642
$ws.Range("I58").Value = 3509.976
$ws.Range("K58").Value = 3509.976
$ws.Range("M58").Value = -3306.976
$ws.Range("H107").Value = 1749
$ws.Range("I107").Value = 1716.6666
$ws.Range("K107").Value = 1716.6666
$ws.Range("M107").Value = 203.3334
$ws.Range("H132").Value = 5452.566
$ws.Range("I132").Value = 1987.4642
$ws.Range("J132").Value = 9333.48
$ws.Range("K132").Value = 5962.392599999999
$ws.Range("L132").Value = 28000.44
$ws.Range("M132").Value = -3432.392599999999
$ws.Range("N132").Value = -33060.44
$ws.Range("H134").Value = 7024.593
$ws.Range("I134").Value = 1410.8667
$ws.Range("K134").Value = 4232.6001
$ws.Range("M134").Value = -1697.6001
$ws.Range("H136").Value = 9516.642
$ws.Range("I136").Value = 3509.976
$ws.Range("K136").Value = 10529.928
$ws.Range("M136").Value = -7979.928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 496
$ws.Range("J23").Value = 114
$ws.Range("L23").Value = 342
$ws.Range("N23").Value = -812
$ws.Range("H87").Value = 18000
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 18800
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 56400
$ws.Range("M87").Value = -28752
$ws.Range("N87").Value = -58896
$ws.Range("H90").Value = 18000
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 18800
$ws.Range("K90").Value = 90000
$ws.Range("L90").Value = 169200
$ws.Range("M90").Value = -83760
$ws.Range("N90").Value = -181680
$ws.Range("H103").Value = 3750
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 3750
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 11250
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -13008
$ws.Range("H107").Value = 7787.643
$ws.Range("I107").Value = 50189.5
$ws.Range("J107").Value = 720.6667
$ws.Range("K107").Value = 150568.5
$ws.Range("L107").Value = 2162.0001
$ws.Range("M107").Value = -148648.5
$ws.Range("N107").Value = -6002.0001
$ws.Range("H114").Value = 2230.75
$ws.Range("I114").Value = 1121.875
$ws.Range("J114").Value = 4448.5
$ws.Range("K114").Value = 3365.625
$ws.Range("L114").Value = 13345.5
$ws.Range("M114").Value = -111.625
$ws.Range("N114").Value = -19853.5
$ws.Range("H115").Value = 1705.2
$ws.Range("I115").Value = 1381.5
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 4144.5
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -2969.5
$ws.Range("N115").Value = -11350
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 7175024
$ws.Range("I122").Value = 20763986
$ws.Range("J122").Value = 1351183.1
$ws.Range("K122").Value = 186875874
$ws.Range("L122").Value = 12160647.9
$ws.Range("M122").Value = -186873424
$ws.Range("N122").Value = -12165547.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 26936.625
$ws.Range("I80").Value = 23664.666
$ws.Range("J80").Value = 28899.8
$ws.Range("K80").Value = 23664.666
$ws.Range("L80").Value = 28899.8
$ws.Range("M80").Value = -22666.666
$ws.Range("N80").Value = -30895.8
$ws.Range("H83").Value = 26936.625
$ws.Range("I83").Value = 23664.666
$ws.Range("J83").Value = 28899.8
$ws.Range("K83").Value = 118323.33
$ws.Range("L83").Value = 144499
$ws.Range("M83").Value = -113331.33
$ws.Range("N83").Value = -154483
$ws.Range("H126").Value = 6077
$ws.Range("I126").Value = 5716.579
$ws.Range("J126").Value = 7218.3335
$ws.Range("K126").Value = 17149.737
$ws.Range("L126").Value = 21655.0005
$ws.Range("M126").Value = -14679.737
$ws.Range("N126").Value = -26595.0005
$ws.Range("H133").Value = 96944.27
$ws.Range("J133").Value = 96944.27
$ws.Range("L133").Value = 96944.27
$ws.Range("N133").Value = -107064.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2566.3572
$ws.Range("I16").Value = 2381.5833
$ws.Range("K16").Value = 2381.5833
$ws.Range("M16").Value = -2211.5833
$ws.Range("H25").Value = 2007
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H132").Value = 10494.926
$ws.Range("I132").Value = 6525.2144
$ws.Range("K132").Value = 19575.6432
$ws.Range("M132").Value = -17045.6432
$ws.Range("H136").Value = 17307.05
$ws.Range("I136").Value = 15964.739
$ws.Range("K136").Value = 47894.217
$ws.Range("M136").Value = -45344.217

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6934.1113
$ws.Range("I132").Value = 3173.9167
$ws.Range("K132").Value = 9521.750100000001
$ws.Range("M132").Value = -6991.750100000001
